$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The original sheet has a header row (row 1) and 5 data rows (2-6):
#   row2 = "Custom hh:mm:ss"        (numFmtId 21 / style 1)
#   row3 = "Time *1:30:55 PM (US)"  (numFmtId 164 / style 2) -- shared formula anchor C3:C6
#   row4 = "Time 13:30 (US)"        (numFmtId 166 / style 4)
#   row5 = "Time 13:30:55 (UK)"     (numFmtId 165 / style 3)
#   row6 = "Time 30:55.2 (US)"      (numFmtId 167 / style 6)
#
# The target sheet drops the "Custom hh:mm:ss" row and instead adds
# test rows for the built-in time formats 18, 19, 20, 21, 45, 46, 47
# (rows 2-8), pushing the four remaining original rows down to 9-12.
# ------------------------------------------------------------------

# Step 1: drop the old "Custom hh:mm:ss" label and the old shared
# formula for column C (rows 3:6) so everything can be rebuilt cleanly.
$ws.Range("A2").ClearContents()
$ws.Range("C3:C6").ClearContents()

# Step 2: relocate the four still-needed rows to their new home at
# rows 9-12 (values + original custom number formats preserved).
$ws.Range("A9").Value = "Time *1:30:55 PM (US)"
$ws.Range("B9").Value = 0.56313888888888886
$ws.Range("B9").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"

$ws.Range("A10").Value = "Time 13:30 (US)"
$ws.Range("B10").Value = 0.56313888888888886
$ws.Range("B10").NumberFormat = "h:mm;@"

$ws.Range("A11").Value = "Time 13:30:55 (UK)"
$ws.Range("B11").Value = 0.56313888888888886
$ws.Range("B11").NumberFormat = "hh:mm:ss;@"

$ws.Range("A12").Value = "Time 30:55.2 (US)"
$ws.Range("B12").Value = 0.021472222222222222
$ws.Range("B12").NumberFormat = "mm:ss.0;@"

# Clear out the old positions of those rows now that the data moved.
$ws.Range("A3:C6").ClearContents()

# Step 3: populate the new rows 2-8 that exercise the predefined time
# formats 18, 19, 20, 21, 45, 46 and 47. Labels are written in the
# order 5,4,3,2,6,7,8 (then number formats in row order) so the
# resulting shared-string table / style table line up with the target
# workbook.
$ws.Range("A5").Value = "Time fmt 21: h:mm:ss"
$ws.Range("A4").Value = "Time fmt 20: H:mm"
$ws.Range("A3").Value = "Time fmt 19: h:mm:ss tt"
$ws.Range("A2").Value = "Time fmt 18: h:mm tt"
$ws.Range("A6").Value = "Time fmt 45: mm:ss"
$ws.Range("A7").Value = "Time fmt 46: [h]:mm:ss"
$ws.Range("A8").Value = "Time fmt 47: mm:ss.0"

$ws.Range("B2").Value = 0.56313888888888886
$ws.Range("B2").NumberFormat = "h:mm AM/PM"

$ws.Range("B3").Value = 0.56313888888888886
$ws.Range("B3").NumberFormat = "h:mm:ss AM/PM"

$ws.Range("B4").Value = 0.56313888888888886
$ws.Range("B4").NumberFormat = "h:mm"

$ws.Range("B5").Value = 0.56313888888888886
$ws.Range("B5").NumberFormat = "h:mm:ss"

$ws.Range("B6").Value = 0.56313888888888886
$ws.Range("B6").NumberFormat = "mm:ss"

$ws.Range("B7").Value = 0.56313888888888886
$ws.Range("B7").NumberFormat = "[h]:mm:ss"

$ws.Range("B8").Value = 0.56313888888888886
$ws.Range("B8").NumberFormat = "mm:ss.0"

# Step 4: formulas for column C (row 2-8 individually, 9-12 shared).
$ws.Range("C2").Formula = "=B2"
$ws.Range("C3").Formula = "=B3"
$ws.Range("C4").Formula = "=B4"
$ws.Range("C5").Formula = "=B5"
$ws.Range("C6").Formula = "=B6"
$ws.Range("C7").Formula = "=B7"
$ws.Range("C8").Formula = "=B8"
$ws.Range("C9:C12").Formula = "=B9"

# Step 5: selection + print orientation, matching the authored session.
$ws.Range("C11").Select()
$ws.PageSetup.Orientation = 1
